# Add a new data row (row 8) to the "Пучиена 2" sheet, mirroring the
# formatting of the row above it (row 7), then fill in the new text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/row look) of the last existing
# data row (A7:E7) down into the new row (A8:E8).
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# Fill in the new row's values. The order below matches the order in
# which the new shared strings were originally authored (English text,
# then filename, then translated text, then converted/mojibake text).
$ws.Cells.Item(8, 3).Value = " It still stinks a little…"
$ws.Cells.Item(8, 1).Value = "SCRIPT/T01P02A/um1104.ssb"
$ws.Cells.Item(8, 4).Value = " Здесь всё ещё пахнет…"
$ws.Cells.Item(8, 5).Value = " Èäåòû âòæ åþæ ðàöîåó…"
$ws.Cells.Item(8, 2).Value = 155

# Match the wrapped-text row height used by the other multi-line rows.
$ws.Rows.Item(8).RowHeight = 43.2

# Leave the selection where the author left it after the edit.
$ws.Range("D6").Select()
